$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab
$ws.Name = "Export as TSV"

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Add error alert title/message to each data validation rule
$validations = @(
    @{ Range = "I2:I1048576"; Title = "Value must come from list"; Message = "Value must be one of: imaging." },
    @{ Range = "J2:J1048576"; Title = "Value must come from list"; Message = "Value must be one of: PAS microscopy." },
    @{ Range = "L2:L1048576"; Title = "Not a boolean"; Message = 'The values in this column must be "TRUE" or "FALSE".' },
    @{ Range = "O2:O1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "P2:P1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." },
    @{ Range = "Q2:Q1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "R2:R1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." },
    @{ Range = "S2:S1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "T2:T1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." }
)

foreach ($item in $validations) {
    $v = $ws.Range($item.Range).Validation
    $v.ErrorTitle = $item.Title
    $v.ErrorMessage = $item.Message
}
